{"js": "// \"Completed Project Charter and Event Table\"\n// Net content changes to the \"Features List\" bullet block:\n//  1. Remove the two early draft bullets \"Generate reports for employees\" and\n//     \"Events can be curated by employees\" (their content is superseded by the\n//     more complete, consistently-worded bullets added in step 2).\n//  2. After \"Members can select if they want to receive promotional emails\",\n//     append the completed set of Employee/Visitor/Member feature bullets,\n//     followed by two blank (but still indented/italic-styled) bullet lines.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1. Delete the two superseded bullets -------------------------------\nconst textsToDelete = [\n  \"Generate reports for employees\",\n  \"Events can be curated by employees\",\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (textsToDelete.indexOf(p.text.trim()) !== -1) {\n    p.delete();\n  }\n}\nawait context.sync();\n\n// --- 2. Insert the new bullets after the \"promotional emails\" line ------\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (\n    paragraphs.items[i].text.trim() ===\n    \"Members can select if they want to receive promotional emails\"\n  ) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (anchor) {\n  const newBullets = [\n    \"Employee can curate available games\",\n    \"Employee can delete a game\",\n    \"Employee can manage events\",\n    \"Employee can view reports\",\n    \"Employee can view a single report\",\n    \"Employee can print a report\",\n    \"Visitor can register account\",\n    \"Member can manage account details\",\n    \"\",\n    \"\",\n  ];\n\n  for (const bulletText of newBullets) {\n    anchor = anchor.insertParagraph(bulletText, Word.InsertLocation.after);\n    // insertParagraph inherits formatting (italic + indent) from the\n    // reference paragraph, matching the style of the rest of the list.\n  }\n}\n\nawait context.sync();\n", "ps1": "# \"Completed Project Charter and Event Table\"\n# Net content changes to the \"Features List\" bullet block:\n#  1. Remove the two early draft bullets \"Generate reports for employees\" and\n#     \"Events can be curated by employees\" (their content is superseded by the\n#     more complete, consistently-worded bullets added in step 2).\n#  2. After \"Members can select if they want to receive promotional emails\",\n#     append the completed set of Employee/Visitor/Member feature bullets,\n#     followed by two blank (but still indented/italic-styled) bullet lines.\n\n$d = $word.ActiveDocument\n\nfunction Delete-ParagraphByText($doc, $text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        if ($p.Range.Text.Trim() -eq $text) {\n            $p.Range.Delete()\n            return $true\n        }\n    }\n    return $false\n}\n\nfunction Find-ParagraphIndexByText($doc, $text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs.Item($i).Range.Text.Trim() -eq $text) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# --- 1. Delete the two superseded bullets --------------------------------\n[void](Delete-ParagraphByText $d \"Generate reports for employees\")\n[void](Delete-ParagraphByText $d \"Events can be curated by employees\")\n\n# --- 2. Insert the new bullets after the \"promotional emails\" line -------\n$anchorIdx = Find-ParagraphIndexByText $d \"Members can select if they want to receive promotional emails\"\n\n$newBullets = @(\n    \"Employee can curate available games\",\n    \"Employee can delete a game\",\n    \"Employee can manage events\",\n    \"Employee can view reports\",\n    \"Employee can view a single report\",\n    \"Employee can print a report\",\n    \"Visitor can register account\",\n    \"Member can manage account details\",\n    \"\",\n    \"\"\n)\n\nif ($anchorIdx -gt 0) {\n    $idx = $anchorIdx\n    foreach ($t in $newBullets) {\n        $cur = $d.Paragraphs.Item($idx)\n        [void]$cur.Range.InsertParagraphAfter()\n        $idx = $idx + 1\n        if ($t -ne \"\") {\n            $newPara = $d.Paragraphs.Item($idx)\n            $newPara.Range.Text = $t\n        }\n    }\n}\n"}
